$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column order: A,B,C,D,E,F,G,H,I,J,K,L,M,N,O,P,Q,R,S,T (20 cols)
$data = New-Object 'object[,]' 15,20

$data[0,0] = "ECs"
$data[0,1] = "Fn1"
$data[0,2] = "Col13a1"
$data[0,3] = "ECs"
$data[0,4] = 3
$data[0,5] = 1
$data[0,6] = 24.721787
$data[0,7] = 74.16536099999999
$data[0,8] = 0.01340847400407589
$data[0,9] = 0.01345828544885583
$data[0,10] = 3
$data[0,11] = 1
$data[0,12] = 0.3015963333333334
$data[0,13] = 0.9047890000000001
$data[0,14] = 0.4737402899861982
$data[0,15] = 0.5687660179545673
$data[0,16] = 7.456000312647665
$data[0,17] = 67.104002813829
$data[0,18] = 0.006352134362963314
$data[0,19] = 0.007654615423241625

$data[1,0] = "ECs"
$data[1,1] = "Fn1"
$data[1,2] = "Col13a1"
$data[1,3] = "FAPs"
$data[1,4] = 3
$data[1,5] = 1
$data[1,6] = 24.721787
$data[1,7] = 74.16536099999999
$data[1,8] = 0.01340847400407589
$data[1,9] = 0.01345828544885583
$data[1,10] = 1
$data[1,11] = 0.3333333333333333
$data[1,12] = 0.01594066666666667
$data[1,13] = 0.047822
$data[1,14] = 0.02503921704145383
$data[1,15] = 0.03006173650500096
$data[1,16] = 0.3940817659713333
$data[1,17] = 3.546735893742
$data[1,18] = 0.0003357376907827479
$data[1,19] = 0.0004045794309725925

$data[2,0] = "ECs"
$data[2,1] = "Fn1"
$data[2,2] = "Col13a1"
$data[2,3] = "sCs"
$data[2,4] = 3
$data[2,5] = 1
$data[2,6] = 24.721787
$data[2,7] = 74.16536099999999
$data[2,8] = 0.01340847400407589
$data[2,9] = 0.01345828544885583
$data[2,10] = 2
$data[2,11] = 1
$data[2,12] = 0.319091
$data[2,13] = 0.638182
$data[2,14] = 0.5012204929723481
$data[2,15] = 0.4011722455404317
$data[2,16] = 7.888499735616999
$data[2,17] = 47.330998413702
$data[2,18] = 0.006720601950329834
$data[2,19] = 0.005399090594641609

$data[3,0] = "FAPs"
$data[3,1] = "Fn1"
$data[3,2] = "Col13a1"
$data[3,3] = "ECs"
$data[3,4] = 3
$data[3,5] = 1
$data[3,6] = 1458.280985666667
$data[3,7] = 4374.842957
$data[3,8] = 0.7909348416823457
$data[3,9] = 0.7938731034993884
$data[3,10] = 3
$data[3,11] = 1
$data[3,12] = 0.3015963333333334
$data[3,13] = 0.9047890000000001
$data[3,14] = 0.4737402899861982
$data[3,15] = 0.5687660179545673
$data[3,16] = 439.8121982467859
$data[3,17] = 3958.309784221073
$data[3,18] = 0.3746977012587822
$data[3,19] = 0.4515280438385812

$data[4,0] = "FAPs"
$data[4,1] = "Fn1"
$data[4,2] = "Col13a1"
$data[4,3] = "FAPs"
$data[4,4] = 3
$data[4,5] = 1
$data[4,6] = 1458.280985666667
$data[4,7] = 4374.842957
$data[4,8] = 0.7909348416823457
$data[4,9] = 0.7938731034993884
$data[4,10] = 1
$data[4,11] = 0.3333333333333333
$data[4,12] = 0.01594066666666667
$data[4,13] = 0.047822
$data[4,14] = 0.02503921704145383
$data[4,15] = 0.03006173650500096
$data[4,16] = 23.24597109885045
$data[4,17] = 209.213739889654
$data[4,18] = 0.01980438916653218
$data[4,19] = 0.02386520405580597

$data[5,0] = "FAPs"
$data[5,1] = "Fn1"
$data[5,2] = "Col13a1"
$data[5,3] = "sCs"
$data[5,4] = 3
$data[5,5] = 1
$data[5,6] = 1458.280985666667
$data[5,7] = 4374.842957
$data[5,8] = 0.7909348416823457
$data[5,9] = 0.7938731034993884
$data[5,10] = 2
$data[5,11] = 1
$data[5,12] = 0.319091
$data[5,13] = 0.638182
$data[5,14] = 0.5012204929723481
$data[5,15] = 0.4011722455404317
$data[5,16] = 465.3243379973623
$data[5,17] = 2791.946027984174
$data[5,18] = 0.3964327512570314
$data[5,19] = 0.3184798556050012

$data[6,0] = "M1"
$data[6,1] = "Fn1"
$data[6,2] = "Col13a1"
$data[6,3] = "ECs"
$data[6,4] = 3
$data[6,5] = 1
$data[6,6] = 60.09979000000001
$data[6,7] = 180.29937
$data[6,8] = 0.03259661091107292
$data[6,9] = 0.03271770480169137
$data[6,10] = 3
$data[6,11] = 1
$data[6,12] = 0.3015963333333334
$data[6,13] = 0.9047890000000001
$data[6,14] = 0.4737402899861982
$data[6,15] = 0.5687660179545673
$data[6,16] = 18.12587629810334
$data[6,17] = 163.13288668293
$data[6,18] = 0.01544232790557895
$data[6,19] = 0.01860871867667102

$data[7,0] = "M1"
$data[7,1] = "Fn1"
$data[7,2] = "Col13a1"
$data[7,3] = "FAPs"
$data[7,4] = 3
$data[7,5] = 1
$data[7,6] = 60.09979000000001
$data[7,7] = 180.29937
$data[7,8] = 0.03259661091107292
$data[7,9] = 0.03271770480169137
$data[7,10] = 1
$data[7,11] = 0.3333333333333333
$data[7,12] = 0.01594066666666667
$data[7,13] = 0.047822
$data[7,14] = 0.02503921704145383
$data[7,15] = 0.03006173650500096
$data[7,16] = 0.9580307191266669
$data[7,17] = 8.622276472140001
$data[7,18] = 0.0008161936154181769
$data[7,19] = 0.0009835510207968506

$data[8,0] = "M1"
$data[8,1] = "Fn1"
$data[8,2] = "Col13a1"
$data[8,3] = "sCs"
$data[8,4] = 3
$data[8,5] = 1
$data[8,6] = 60.09979000000001
$data[8,7] = 180.29937
$data[8,8] = 0.03259661091107292
$data[8,9] = 0.03271770480169137
$data[8,10] = 2
$data[8,11] = 1
$data[8,12] = 0.319091
$data[8,13] = 0.638182
$data[8,14] = 0.5012204929723481
$data[8,15] = 0.4011722455404317
$data[8,16] = 19.17730209089
$data[8,17] = 115.06381254534
$data[8,18] = 0.01633808939007579
$data[8,19] = 0.01312543510422349

$data[9,0] = "M2"
$data[9,1] = "Fn1"
$data[9,2] = "Col13a1"
$data[9,3] = "ECs"
$data[9,4] = 3
$data[9,5] = 1
$data[9,6] = 280.168911
$data[9,7] = 840.5067330000001
$data[9,8] = 0.1519565539454633
$data[9,9] = 0.152521060800867
$data[9,10] = 3
$data[9,11] = 1
$data[9,12] = 0.3015963333333334
$data[9,13] = 0.9047890000000001
$data[9,14] = 0.4737402899861982
$data[9,15] = 0.5687660179545673
$data[9,16] = 84.49791627159301
$data[9,17] = 760.4812464443371
$data[9,18] = 0.07198794193142717
$data[9,19] = 0.08674879640591558

$data[10,0] = "M2"
$data[10,1] = "Fn1"
$data[10,2] = "Col13a1"
$data[10,3] = "FAPs"
$data[10,4] = 3
$data[10,5] = 1
$data[10,6] = 280.168911
$data[10,7] = 840.5067330000001
$data[10,8] = 0.1519565539454633
$data[10,9] = 0.152521060800867
$data[10,10] = 1
$data[10,11] = 0.3333333333333333
$data[10,12] = 0.01594066666666667
$data[10,13] = 0.047822
$data[10,14] = 0.02503921704145383
$data[10,15] = 0.03006173650500096
$data[10,16] = 4.466079220614001
$data[10,17] = 40.194712985526
$data[10,18] = 0.003804873135111844
$data[10,19] = 0.004585047941258896

$data[11,0] = "M2"
$data[11,1] = "Fn1"
$data[11,2] = "Col13a1"
$data[11,3] = "sCs"
$data[11,4] = 3
$data[11,5] = 1
$data[11,6] = 280.168911
$data[11,7] = 840.5067330000001
$data[11,8] = 0.1519565539454633
$data[11,9] = 0.152521060800867
$data[11,10] = 2
$data[11,11] = 1
$data[11,12] = 0.319091
$data[11,13] = 0.638182
$data[11,14] = 0.5012204929723481
$data[11,15] = 0.4011722455404317
$data[11,16] = 89.39937797990102
$data[11,17] = 536.3962678794061
$data[11,18] = 0.07616373887892433
$data[11,19] = 0.06118721645369254

$data[12,0] = "sCs"
$data[12,1] = "Fn1"
$data[12,2] = "Col13a1"
$data[12,3] = "ECs"
$data[12,4] = 2
$data[12,5] = 1
$data[12,6] = 20.472042
$data[12,7] = 40.944084
$data[12,8] = 0.01110351945704208
$data[12,9] = 0.00742984544919738
$data[12,10] = 3
$data[12,11] = 1
$data[12,12] = 0.3015963333333334
$data[12,13] = 0.9047890000000001
$data[12,14] = 0.4737402899861982
$data[12,15] = 0.5687660179545673
$data[12,16] = 6.174292803046001
$data[12,17] = 37.045756818276
$data[12,18] = 0.005260184527446509
$data[12,19] = 0.004225843610157856

$data[13,0] = "sCs"
$data[13,1] = "Fn1"
$data[13,2] = "Col13a1"
$data[13,3] = "FAPs"
$data[13,4] = 2
$data[13,5] = 1
$data[13,6] = 20.472042
$data[13,7] = 40.944084
$data[13,8] = 0.01110351945704208
$data[13,9] = 0.00742984544919738
$data[13,10] = 1
$data[13,11] = 0.3333333333333333
$data[13,12] = 0.01594066666666667
$data[13,13] = 0.047822
$data[13,14] = 0.02503921704145383
$data[13,15] = 0.03006173650500096
$data[13,16] = 0.3263379975080001
$data[13,17] = 1.958027985048
$data[13,18] = 0.0002780234336088823
$data[13,19] = 0.0002233540561666522

$data[14,0] = "sCs"
$data[14,1] = "Fn1"
$data[14,2] = "Col13a1"
$data[14,3] = "sCs"
$data[14,4] = 2
$data[14,5] = 1
$data[14,6] = 20.472042
$data[14,7] = 40.944084
$data[14,8] = 0.01110351945704208
$data[14,9] = 0.00742984544919738
$data[14,10] = 2
$data[14,11] = 1
$data[14,12] = 0.319091
$data[14,13] = 0.638182
$data[14,14] = 0.5012204929723481
$data[14,15] = 0.4011722455404317
$data[14,16] = 6.532444353822001
$data[14,17] = 26.129777415288
$data[14,18] = 0.00556531149598669
$data[14,19] = 0.00298064778287287

$ws.Range("A2:T16").Value = $data
